# Update the "time_taken" timestamps (column F) on the "data" sheet to the
# new values recorded when the metadata snapshot was refreshed.
$newTimes = @(
    "2021-10-05 14:33:57.597292",
    "2021-10-05 14:33:57.597299",
    "2021-10-05 14:33:57.597302",
    "2021-10-05 14:33:57.597305",
    "2021-10-05 14:33:57.597308",
    "2021-10-05 14:33:57.597310",
    "2021-10-05 14:33:57.597313",
    "2021-10-05 14:33:57.597315",
    "2021-10-05 14:33:57.597318",
    "2021-10-05 14:33:57.597321",
    "2021-10-05 14:33:57.597323",
    "2021-10-05 14:33:57.597326",
    "2021-10-05 14:33:57.597328",
    "2021-10-05 14:33:57.597331",
    "2021-10-05 14:33:57.597333",
    "2021-10-05 14:33:57.597336",
    "2021-10-05 14:33:57.597338",
    "2021-10-05 14:33:57.597341",
    "2021-10-05 14:33:57.597344",
    "2021-10-05 14:33:57.597346",
    "2021-10-05 14:33:57.597349",
    "2021-10-05 14:33:57.597351",
    "2021-10-05 14:33:57.597354",
    "2021-10-05 14:33:57.597356",
    "2021-10-05 14:33:57.597359",
    "2021-10-05 14:33:57.597361",
    "2021-10-05 14:33:57.597364",
    "2021-10-05 14:33:57.597366",
    "2021-10-05 14:33:57.597369",
    "2021-10-05 14:33:57.597371",
    "2021-10-05 14:33:57.597374",
    "2021-10-05 14:33:57.597376",
    "2021-10-05 14:33:57.597379",
    "2021-10-05 14:33:57.597382",
    "2021-10-05 14:33:57.597385",
    "2021-10-05 14:33:57.597387",
    "2021-10-05 14:33:57.597390",
    "2021-10-05 14:33:57.597392",
    "2021-10-05 14:33:57.597395",
    "2021-10-05 14:33:57.597397",
    "2021-10-05 14:33:57.597400",
    "2021-10-05 14:33:57.597403",
    "2021-10-05 14:33:57.597405",
    "2021-10-05 14:33:57.597408",
    "2021-10-05 14:33:57.597410",
    "2021-10-05 14:33:57.597413",
    "2021-10-05 14:33:57.597415",
    "2021-10-05 14:33:57.597418",
    "2021-10-05 14:33:57.597420",
    "2021-10-05 14:33:57.597423",
    "2021-10-05 14:33:57.597425",
    "2021-10-05 14:33:57.597428",
    "2021-10-05 14:33:57.597431",
    "2021-10-05 14:33:57.597433",
    "2021-10-05 14:33:57.597436",
    "2021-10-05 14:33:57.597438",
    "2021-10-05 14:33:57.597441",
    "2021-10-05 14:33:57.597443",
    "2021-10-05 14:33:57.597446",
    "2021-10-05 14:33:57.597448",
    "2021-10-05 14:33:57.597451",
    "2021-10-05 14:33:57.597453",
    "2021-10-05 14:33:57.597455",
    "2021-10-05 14:33:57.597458",
    "2021-10-05 14:33:57.597461",
    "2021-10-05 14:33:57.597464",
    "2021-10-05 14:33:57.597467",
    "2021-10-05 14:33:57.597469",
    "2021-10-05 14:33:57.597472",
    "2021-10-05 14:33:57.597474",
    "2021-10-05 14:33:57.597476",
    "2021-10-05 14:33:57.597479",
    "2021-10-05 14:33:57.597481",
    "2021-10-05 14:33:57.597484",
    "2021-10-05 14:33:57.597486",
    "2021-10-05 14:33:57.597489",
    "2021-10-05 14:33:57.597493",
    "2021-10-05 14:33:57.597496",
    "2021-10-05 14:33:57.597498",
    "2021-10-05 14:33:57.597501",
    "2021-10-05 14:33:57.597503",
    "2021-10-05 14:33:57.597506",
    "2021-10-05 14:33:57.597508",
    "2021-10-05 14:33:57.597511",
    "2021-10-05 14:33:57.597513",
    "2021-10-05 14:33:57.597516",
    "2021-10-05 14:33:57.597518",
    "2021-10-05 14:33:57.597521",
    "2021-10-05 14:33:57.597523",
    "2021-10-05 14:33:57.597526",
    "2021-10-05 14:33:57.597528",
    "2021-10-05 14:33:57.597530",
    "2021-10-05 14:33:57.597564",
    "2021-10-05 14:33:57.597569",
    "2021-10-05 14:33:57.597572",
    "2021-10-05 14:33:57.597575",
    "2021-10-05 14:33:57.597577",
    "2021-10-05 14:33:57.597580",
    "2021-10-05 14:33:57.597582",
    "2021-10-05 14:33:57.597585",
    "2021-10-05 14:33:57.597588",
    "2021-10-05 14:33:57.597590",
    "2021-10-05 14:33:57.597593",
    "2021-10-05 14:33:57.597595",
    "2021-10-05 14:33:57.597598",
    "2021-10-05 14:33:57.597600",
    "2021-10-05 14:33:57.597603"
)

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

for ($i = 0; $i -lt $newTimes.Count; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# Add a new "metadata" worksheet after the "data" sheet, describing the
# panelapp source this data snapshot was pulled from.
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Hereditary Neuropathy_CMT - isolated"
$metaSheet.Range("C2").Value = 3069
$metaSheet.Range("D2").Value = "'1.9"
$metaSheet.Range("D2").Style = "Normal"
$metaSheet.Range("E2").Value = "2021-10-04T06:47:47.252923Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:33:57.593992"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3069/?format=json"

# Match the bold/centered/bordered header style used on the "data" sheet
# (style index reused, not re-created) for the header row and index cell.
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Keep "data" as the active sheet/selection, matching the unchanged
# <bookViews> (activeTab stays 0) in the source diff.
$dataSheet.Activate()
$dataSheet.Range("A1").Select()

Write-Output "done"
